# Updates the Cambodia CPL 2023-2024 sheet:
#  1. Rows 59 and 60 had been entered with swapped match data; fix by
#     swapping the F:V (home..url) contents between the two rows.
#  2. Two newly played matches are appended as rows 77 and 78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the match data (columns F..V) between rows 59 and 60 ---------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $cols) {
    $addrOld = $col + "59"
    $addrNew = $col + "60"
    $valOld = $ws.Range($addrOld).Value2
    $valNew = $ws.Range($addrNew).Value2
    $ws.Range($addrOld).Value2 = $valNew
    $ws.Range($addrNew).Value2 = $valOld
}

# --- 2. Append the two new match rows (77 and 78) -------------------------
# Copy formatting (styles) from the last existing data row so the new rows
# match the rest of the table (bold/bordered index column, date number
# format on column E, etc.)
$ws.Range("A76:V76").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$ws.Range("A76:V76").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)

# Row 77
$ws.Range("A77").Value2 = 76
$ws.Range("B77").Value2 = "cambodia"
$ws.Range("C77").Value2 = "cpl"
$ws.Range("D77").Value2 = "2023-2024"
$ws.Range("E77").Value2 = 45297.5
$ws.Range("F77").Value2 = "Dangkor"
$ws.Range("G77").Value2 = 1
$ws.Range("H77").Value2 = "Kirivong Sok Sen Chey"
$ws.Range("I77").Value2 = 1
$ws.Range("J77").Value2 = 1.98
$ws.Range("K77").Value2 = "06/01/2024 01:12"
$ws.Range("L77").Value2 = 1.73
$ws.Range("M77").Value2 = "06/01/2024 11:59"
$ws.Range("N77").Value2 = 3.56
$ws.Range("O77").Value2 = "06/01/2024 01:12"
$ws.Range("P77").Value2 = 3.76
$ws.Range("Q77").Value2 = "06/01/2024 11:59"
$ws.Range("R77").Value2 = 3.05
$ws.Range("S77").Value2 = "06/01/2024 01:12"
$ws.Range("T77").Value2 = 3.82
$ws.Range("U77").Value2 = "06/01/2024 11:59"
$ws.Range("V77").Value2 = "https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-kirivong-sok-sen-chey/Mw7kIejg/"

# Row 78
$ws.Range("A78").Value2 = 77
$ws.Range("B78").Value2 = "cambodia"
$ws.Range("C78").Value2 = "cpl"
$ws.Range("D78").Value2 = "2023-2024"
$ws.Range("E78").Value2 = 45297.5
$ws.Range("F78").Value2 = "NagaWorld"
$ws.Range("G78").Value2 = 0
$ws.Range("H78").Value2 = "Boeung Ket"
$ws.Range("I78").Value2 = 0
$ws.Range("J78").Value2 = 3.42
$ws.Range("K78").Value2 = "06/01/2024 01:12"
$ws.Range("L78").Value2 = 3.47
$ws.Range("M78").Value2 = "06/01/2024 11:44"
$ws.Range("N78").Value2 = 3.71
$ws.Range("O78").Value2 = "06/01/2024 01:12"
$ws.Range("P78").Value2 = 3.74
$ws.Range("Q78").Value2 = "06/01/2024 10:05"
$ws.Range("R78").Value2 = 1.81
$ws.Range("S78").Value2 = "06/01/2024 01:12"
$ws.Range("T78").Value2 = 1.82
$ws.Range("U78").Value2 = "06/01/2024 11:44"
$ws.Range("V78").Value2 = "https://www.betexplorer.com/football/cambodia/cpl/nagaworld-boeung-ket/G06gHF5a/"
